$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 123.67
$ws.Range("I15").Value = 123.67
$ws.Range("K15").Value = 371.01
$ws.Range("M15").Value = -202.01
$ws.Range("H32").Value = 16667697
$ws.Range("I32").Value = 83333850
$ws.Range("J32").Value = 1159
$ws.Range("K32").Value = 83333850
$ws.Range("L32").Value = 1159
$ws.Range("M32").Value = -83333524
$ws.Range("N32").Value = -1811
$ws.Range("H96").Value = 62569970
$ws.Range("I96").Value = 5161.857
$ws.Range("J96").Value = 111231490
$ws.Range("K96").Value = 15485.571
$ws.Range("L96").Value = 333694470
$ws.Range("M96").Value = -14112.571
$ws.Range("N96").Value = -333697216
$ws.Range("H100").Value = 2892.3076
$ws.Range("I100").Value = 3260
$ws.Range("J100").Value = 2662.5
$ws.Range("K100").Value = 3260
$ws.Range("L100").Value = 2662.5
$ws.Range("M100").Value = -2719
$ws.Range("N100").Value = -3744.5
$ws.Range("H109").Value = 43684
$ws.Range("J109").Value = 43684
$ws.Range("L109").Value = 43684
$ws.Range("N109").Value = -46458
$ws.Range("H114").Value = 41212
$ws.Range("J114").Value = 41212
$ws.Range("L114").Value = 41212
$ws.Range("N114").Value = -49890
$ws.Range("H133").Value = 38627.367
$ws.Range("J133").Value = 38627.367
$ws.Range("L133").Value = 38627.367
$ws.Range("N133").Value = -48747.367
$ws.Range("H137").Value = 4008.5305
$ws.Range("I137").Value = 1901
$ws.Range("J137").Value = 4098.213
$ws.Range("K137").Value = 5703
$ws.Range("L137").Value = 12294.639
$ws.Range("M137").Value = -3153
$ws.Range("N137").Value = -17394.639

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3076.7234
$ws.Range("I2").Value = 3897.9143
$ws.Range("J2").Value = 681.5833
$ws.Range("K2").Value = 3897.9143
$ws.Range("L2").Value = 681.5833
$ws.Range("M2").Value = -3784.9143
$ws.Range("N2").Value = -907.5833
$ws.Range("H32").Value = 32945.24
$ws.Range("I32").Value = 32170.29
$ws.Range("J32").Value = 38702
$ws.Range("K32").Value = 32170.29
$ws.Range("L32").Value = 38702
$ws.Range("M32").Value = -31883.29
$ws.Range("N32").Value = -39276
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H61").Value = 3049.932
$ws.Range("I61").Value = 1490.1666
$ws.Range("J61").Value = 4129.769
$ws.Range("K61").Value = 1490.1666
$ws.Range("L61").Value = 4129.769
$ws.Range("M61").Value = -1278.1666
$ws.Range("N61").Value = -4553.769
$ws.Range("H107").Value = 35004.75
$ws.Range("J107").Value = 35004.75
$ws.Range("L107").Value = 35004.75
$ws.Range("N107").Value = -42684.75
$ws.Range("H109").Value = 49877
$ws.Range("J109").Value = 49877
$ws.Range("L109").Value = 49877
$ws.Range("N109").Value = -52651
$ws.Range("H116").Value = 3076.7234
$ws.Range("I116").Value = 3897.9143
$ws.Range("J116").Value = 681.5833
$ws.Range("K116").Value = 3897.9143
$ws.Range("L116").Value = 681.5833
$ws.Range("M116").Value = -1603.9143
$ws.Range("N116").Value = -5269.5833
$ws.Range("H118").Value = 49997.332
$ws.Range("J118").Value = 49997.332
$ws.Range("L118").Value = 49997.332
$ws.Range("N118").Value = -53311.332
$ws.Range("H131").Value = 50399.5
$ws.Range("J131").Value = 50399.5
$ws.Range("L131").Value = 50399.5
$ws.Range("N131").Value = -60479.5
$ws.Range("H132").Value = 13890714
$ws.Range("I132").Value = 23810610
$ws.Range("J132").Value = 2857.8667
$ws.Range("K132").Value = 71431830
$ws.Range("L132").Value = 8573.6001
$ws.Range("M132").Value = -71429300
$ws.Range("N132").Value = -13633.6001
$ws.Range("H136").Value = 3049.932
$ws.Range("I136").Value = 1490.1666
$ws.Range("J136").Value = 4129.769
$ws.Range("K136").Value = 4470.4998
$ws.Range("L136").Value = 12389.307
$ws.Range("M136").Value = -1920.4998
$ws.Range("N136").Value = -17489.307

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3076.7234
$ws.Range("I3").Value = 3897.9143
$ws.Range("J3").Value = 681.5833
$ws.Range("K3").Value = 3897.9143
$ws.Range("L3").Value = 681.5833
$ws.Range("M3").Value = -3783.9143
$ws.Range("N3").Value = -909.5833
$ws.Range("H35").Value = 31590
$ws.Range("J35").Value = 31590
$ws.Range("L35").Value = 31590
$ws.Range("N35").Value = -32210
$ws.Range("H100").Value = 41996
$ws.Range("J100").Value = 41996
$ws.Range("L100").Value = 41996
$ws.Range("N100").Value = -44160
$ws.Range("H108").Value = 37621
$ws.Range("J108").Value = 37621
$ws.Range("L108").Value = 37621
$ws.Range("N108").Value = -45301

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2376
$ws.Range("I22").Value = 394.2857
$ws.Range("K22").Value = 394.2857
$ws.Range("M22").Value = -44.28570000000002
$ws.Range("H110").Value = 38003.332
$ws.Range("J110").Value = 38003.332
$ws.Range("L110").Value = 38003.332
$ws.Range("N110").Value = -46183.332
$ws.Range("H111").Value = 40166.8
$ws.Range("J111").Value = 40166.8
$ws.Range("L111").Value = 40166.8
$ws.Range("N111").Value = -48346.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1646.1875
$ws.Range("J34").Value = 1742.6
$ws.Range("L34").Value = 5227.799999999999
$ws.Range("N34").Value = -5395.799999999999
$ws.Range("H39").Value = 2300.4443
$ws.Range("J39").Value = 2672
$ws.Range("L39").Value = 8016
$ws.Range("N39").Value = -8604
$ws.Range("H55").Value = 1580.8
$ws.Range("I55").Value = 704
$ws.Range("J55").Value = 1800
$ws.Range("K55").Value = 2112
$ws.Range("L55").Value = 5400
$ws.Range("M55").Value = -1935
$ws.Range("N55").Value = -5754
$ws.Range("H75").Value = 933.3333
$ws.Range("I75").Value = 933.3333
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2799.9999
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -1801.9999
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 933.3333
$ws.Range("I78").Value = 933.3333
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 8399.9997
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -3407.9997
$ws.Range("N78").ClearContents()
$ws.Range("H92").Value = 1144.9445
$ws.Range("I92").Value = 1064
$ws.Range("J92").Value = 1272.1428
$ws.Range("K92").Value = 3192
$ws.Range("L92").Value = 3816.4284
$ws.Range("M92").Value = -1944
$ws.Range("N92").Value = -6312.428400000001
$ws.Range("H132").Value = 1774.5
$ws.Range("I132").Value = 1100.5
$ws.Range("J132").Value = 1999.1666
$ws.Range("K132").Value = 9904.5
$ws.Range("L132").Value = 17992.4994
$ws.Range("M132").Value = -7374.5
$ws.Range("N132").Value = -23052.4994

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 49734
$ws.Range("J116").Value = 49734
$ws.Range("L116").Value = 49734
$ws.Range("N116").Value = -58912
$ws.Range("H126").Value = 6277.5
$ws.Range("I126").Value = 7912.8125
$ws.Range("K126").Value = 23738.4375
$ws.Range("M126").Value = -21268.4375

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H93").Value = 1725
$ws.Range("I93").Value = 1033.3334
$ws.Range("J93").Value = 3800
$ws.Range("K93").Value = 1033.3334
$ws.Range("L93").Value = 3800
$ws.Range("M93").Value = 214.6666
$ws.Range("N93").Value = -6296
$ws.Range("H108").Value = 48622
$ws.Range("J108").Value = 48622
$ws.Range("L108").Value = 48622
$ws.Range("N108").Value = -56302
$ws.Range("H109").Value = 29324.5
$ws.Range("J109").Value = 29324.5
$ws.Range("L109").Value = 29324.5
$ws.Range("N109").Value = -32098.5
$ws.Range("H111").Value = 41171.4
$ws.Range("J111").Value = 41171.4
$ws.Range("L111").Value = 41171.4
$ws.Range("N111").Value = -49351.4
$ws.Range("H127").Value = 50557.668
$ws.Range("J127").Value = 50557.668
$ws.Range("L127").Value = 50557.668
$ws.Range("N127").Value = -60477.668
$ws.Range("H132").Value = 4893.231
$ws.Range("I132").Value = 7371.7144
$ws.Range("J132").Value = 3980.1052
$ws.Range("K132").Value = 22115.1432
$ws.Range("L132").Value = 11940.3156
$ws.Range("M132").Value = -19585.1432
$ws.Range("N132").Value = -17000.3156

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 44750.5
$ws.Range("J16").Value = 44750.5
$ws.Range("L16").Value = 44750.5
$ws.Range("N16").Value = -45334.5
$ws.Range("H40").Value = 12081.667
$ws.Range("J40").Value = 12798
$ws.Range("L40").Value = 12798
$ws.Range("N40").Value = -13096
$ws.Range("H81").Value = 1724.2858
$ws.Range("J81").Value = 1800
$ws.Range("L81").Value = 3600
$ws.Range("N81").Value = -5722
$ws.Range("H84").Value = 1724.2858
$ws.Range("J84").Value = 1800
$ws.Range("L84").Value = 18000
$ws.Range("N84").Value = -28608
$ws.Range("H110").Value = 46974.668
$ws.Range("J110").Value = 46974.668
$ws.Range("L110").Value = 46974.668
$ws.Range("N110").Value = -55154.668
$ws.Range("H126").Value = 1089702.5
$ws.Range("I126").Value = 1225852.5
$ws.Range("J126").Value = 503.33334
$ws.Range("K126").Value = 3677557.5
$ws.Range("L126").Value = 1510.00002
$ws.Range("M126").Value = -3675087.5
$ws.Range("N126").Value = -6450.000019999999
$ws.Range("H132").Value = 2020.5454
$ws.Range("I132").Value = 1355.8096
$ws.Range("J132").Value = 3183.8333
$ws.Range("K132").Value = 4067.4288
$ws.Range("L132").Value = 9551.499899999999
$ws.Range("M132").Value = -1537.4288
$ws.Range("N132").Value = -14611.4999
$ws.Range("H137").Value = 47178.5
$ws.Range("J137").Value = 47178.5
$ws.Range("L137").Value = 47178.5
$ws.Range("N137").Value = -57378.5
